# Atualização do select cimento
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 4.8
$ws.Range("B3").Value = 18.24
$ws.Range("B4").Value = 72
$ws.Range("B5").Value = 456
$ws.Range("B6").Value = 10944
$ws.Range("B7").Value = 8
$ws.Range("B8").Value = 16
$ws.Range("B9").Value = 8
